# "last bit o cleanup"
#
# 1. On slide 1, flesh out the "Architectural:" bullet list in the content
#    placeholder: split "Collapse Business Services" off from its old
#    "(DNSAdapter+DNSServiceImpl)" qualifier and add three new sub-bullets
#    underneath it, rename "datastore" -> "data store", and re-split the
#    last two one-liners ("PTR Record Type Implementation" /
#    "Domain Registration Pass-through" -> "Domain Registration").
# 2. Drop the second slide ("POSSIBLE FUTURE FEATURES (cont'd)") entirely --
#    it only ever held two empty placeholder paragraphs.

$p = $ppt.ActivePresentation
$CR = [string][char]13

# --- Slide 1: rewrite the content placeholder bullets -----------------------
$s1 = $p.Slides.Item(1)
$body = $s1.Shapes.Item(2).TextFrame.TextRange

# Para 2: "Collapse Business Services (DNSAdapter+DNSServiceImpl)"
#      -> "Collapse Business " + "Services"
$para = $body.Paragraphs(2, 1)
$para.Text = "Collapse Business "
[void]$para.InsertAfter("Services")

# Insert three new lvl-3 sub-bullets right after it.
[void]$para.InsertAfter($CR)
$para = $body.Paragraphs(3, 1)
$para.Text = "Reduce Adapter to simple SPI with no code"
$para.IndentLevel = 3

[void]$para.InsertAfter($CR)
$para = $body.Paragraphs(4, 1)
$para.Text = "Consolidate code into service provider(s)"
$para.IndentLevel = 3

[void]$para.InsertAfter($CR)
$para = $body.Paragraphs(5, 1)
$para.Text = "Better leverage OSGI service registry"
$para.IndentLevel = 3

# Para (now 6): "Replace MySQL datastore (with say Redis)"
#      -> "Replace MySQL " + "data store " + "(with say " + "Redis" + ")"
$para = $body.Paragraphs(6, 1)
$para.Text = "Replace MySQL "
[void]$para.InsertAfter("data store ")
[void]$para.InsertAfter("(with say ")
[void]$para.InsertAfter("Redis")
[void]$para.InsertAfter(")")

# Para (now 7): "PTR Record Type Implementation" -> "PTR " + "Record Type Implementation"
$para = $body.Paragraphs(7, 1)
$para.Text = "PTR "
[void]$para.InsertAfter("Record Type Implementation")

# Para (now 8): "Domain Registration Pass-through" -> "Domain " + "Registration"
$para = $body.Paragraphs(8, 1)
$para.Text = "Domain "
[void]$para.InsertAfter("Registration")

# The placeholder used to end with two identical empty "no bullet" paragraphs;
# only one of them survives.
$body.Paragraphs(10, 1).Delete()

# --- Drop the second slide (and its notes page) entirely --------------------
$p.Slides.Item(2).Delete()
